# Applies the "cleaned up input and output" edit to Max's Budget.xlsx
# Re-sorts the Income (A/B) and Expense (C/D) category lists alphabetically
# and updates the budget amounts next to each category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Income categories (columns A & B), now sorted alphabetically ---
$income = @(
    @("Bonuses", 0),
    @("Interest Income", 2000),
    @("Investments", 0),
    @("Paychecks", 3000),
    @("Reimbursements", 0),
    @("Rental Incomes", 1000),
    @("Returned Purchases", 0)
)

for ($i = 0; $i -lt $income.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $income[$i][0]
    $ws.Cells.Item($row, 2).Value = $income[$i][1]
}

# --- Expense categories (columns C & D), now sorted alphabetically ---
$expenses = @(
    @("Auto Repair/Transportation", 0),
    @("Clothing", 300),
    @("Debt & Interest Payments", 30),
    @("Eating Out/Delivery", 10),
    @("Education Expenses", 300),
    @("Electronics/Virtual Products", 40),
    @("Fees & Charges", 500),
    @("Gas", 10),
    @("Gifts/Donations", 20),
    @("Groceries", 30),
    @("Health/Medical", 40),
    @("Movies", 50),
    @("Music", 0),
    @("Personal Care", 60),
    @("Pets", 0),
    @("Phone", 10),
    @("Rent", 20),
    @("Streaming Services/Subscriptions", 0),
    @("Travel", 300),
    @("Utilities", 400)
)

for ($i = 0; $i -lt $expenses.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $expenses[$i][0]
    $ws.Cells.Item($row, 4).Value = $expenses[$i][1]
}
